$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B13").Value = "Son"
$ws.Range("C13").Value = "Romain"

$ws.Range("B14").Value = "Graphisme"
$ws.Range("C14").Value = "Steve"

$ws.Range("B15").Value = "Capacités spéciales (Decorator)"
$ws.Range("B16").Value = "Déplacement de blocs"
$ws.Range("B17").Value = "Menu pause amélioré"
$ws.Range("B18").Value = "Eviter les messages infos du tuto"
$ws.Range("B19").Value = "Interaction avec les acteurs (loutre)"
$ws.Range("B20").Value = "Crédits du jeu"

$ws.Range("E15").Select()
